# Weekly data refresh: insert one new price record for "Poroto verde" at
# Vega Modelo de Temuco (row 125), pushing the existing historical rows
# (125-205) down by one (to 126-206).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 125; rows 125..205 shift down to 126..206.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(125, 1).Value  = 10
$ws.Cells.Item(125, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(125, 3).Value  = "La Araucanía"
$ws.Cells.Item(125, 4).Value  = 45029
$ws.Cells.Item(125, 5).Value  = 9
$ws.Cells.Item(125, 6).Value  = 100112031
$ws.Cells.Item(125, 7).Value  = "Poroto verde"
$ws.Cells.Item(125, 8).Value  = "Sin especificar"
$ws.Cells.Item(125, 9).Value  = "Primera"
$ws.Cells.Item(125, 10).Value = 40
$ws.Cells.Item(125, 11).Value = 25000
$ws.Cells.Item(125, 12).Value = 25000
$ws.Cells.Item(125, 13).Value = 25000
$ws.Cells.Item(125, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(125, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(125, 16).Value = 1000
$ws.Cells.Item(125, 17).Value = 25
$ws.Cells.Item(125, 18).Value = "Hortaliza"
